# Edited UI sequence diagram in DG
# Reposition / resize three shapes on the (single) slide of the
# LogicComponentSequenceDiagram deck, matching the author's manual
# diagram tidy-up.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rectangle 5" (activation bar) - made narrower and taller, same top-left.
$rect = $s.Shapes.Item("Rectangle 5")
$rect.Left   = 118.17724609375
$rect.Top    = 237.07962036132812
$rect.Width  = 11.052835464477539
$rect.Height = 266.92041015625

# "Straight Arrow Connector 36" - moved down and slightly left.
$conn36 = $s.Shapes.Item("Straight Arrow Connector 36")
$conn36.Left = 22.92259979248047
$conn36.Top  = 498.0

# "Straight Arrow Connector 76" - moved straight down.
$conn76 = $s.Shapes.Item("Straight Arrow Connector 76")
$conn76.Top = 486.34576416015625
